$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that follows the H1 title.
#    (<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#     <w:r><w:t>: Discover African Quest ...</w:t></w:r></w:p>)
# ---------------------------------------------------------------------
$metaLabel = "Meta description"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith($metaLabel)) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Replace the DALLE image-prompt paragraph's text with the meta
#    description copy (keeps the paragraph's existing italic run).
# ---------------------------------------------------------------------
$oldText = 'Hello DALLE, I would like you to create a feature image for the slot game "African Quest". The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be standing in front of an African savannah, surrounded by African animals such as elephants, giraffes, zebras, and rhinos. The warrior should be holding a map of Africa in one hand and a bag of gold coins in the other hand. The background color of the image should be orange or yellow to evoke the warmth of the African sun. Please make sure that the image is eye-catching and appealing to slot game players. Thank you!'
$descText = "Discover African Quest slot game with special features, high volatility, and RTP. Play for free and win up to €84,000 with this African-themed game."
[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $descText, 2)

# ---------------------------------------------------------------------
# 3) Insert a new bold paragraph ("Play African Quest Free: Review of
#    Gameplay, Design & Winnings") right before that paragraph.
# ---------------------------------------------------------------------
$descStart = "Discover African Quest slot game"
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith($descStart)) {
        $targetIdx = $i
        break
    }
}

$targetPara = $d.Paragraphs($targetIdx)
[void]$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($targetIdx)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play African Quest Free: Review of Gameplay, Design &amp; Winnings</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($titleXml)
